$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.095.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "'2.502.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'320.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "'107.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("D7").Value = "'0.526"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").Value = "'39.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").Value = "'20.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.13%  "

$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").Value = "'7.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").Value = "'2.894.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").Value = "'2.503.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").Value = "'47.962.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").Value = "'12.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("D22").Value = "'2.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").Value = "'276.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.49%  "

$ws.Range("D24").Value = "'71.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").Value = "'25.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "'0.142"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("D30").Value = "'9.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "

$ws.Range("D31").Value = "'34.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("D33").Value = "'19.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").Value = "'0.0782"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").Value = "'4.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("D39").Value = "'2.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "'21.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").Value = "'0.0302"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").Value = "'2.023.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("E46").Value = "  +2.27%  "

$ws.Range("E47").Value = "  -1.59%  "

$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("D49").Value = "'8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("D51").Value = "'80.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.38%  "
